$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F (shifts old F:J -> G:K), copying the
# formatting of the column to its left (E), which is what Excel does.
$ws.Columns("F").Insert()

# The new "ao_ref_id" column needs to look like the other field-code
# columns (G:K) rather than inheriting column E's style, so copy the
# formatting from G4 down into row 4 of the new column.
$ws.Range("G4").Copy() | Out-Null
$ws.Range("F4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# New column F content (Archival Object Ref ID / ao_ref_id)
$ws.Range("F2").Value = "Archival Object Basic Information"
$ws.Range("F3").Value = "Ref_id"
$ws.Range("F4").Value = "ao_ref_id"
$ws.Range("F5").Value = "REF ID"

# Give the new column its own (narrower) width
$ws.Columns("F").ColumnWidth = 9.86

# Fix up the existing C3 comment text (was "EAD ID REQUIRED IF NO URI")
$ws.Range("C3").Value = "EAD ID "

# Restore the last active selection to G1
$ws.Range("G1").Select() | Out-Null
